$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("E5").Value = 1540
$ws.Range("H5").Value = 45141.04206952546
$ws.Range("J5").Value = '07/06/23 08:49'
$ws.Range("K5").Value = '07/06/23 08:49'
$ws.Range("M5").Value = '$1,540 as of 7/5/2023 11:37:54 AM'
$ws.Range("N5").Value = 1580

# Row 7
$ws.Range("A7").Value = 'L682801'
$ws.Range("C7").Value = 'SB#5'
$ws.Range("E7").Value = 2100
$ws.Range("H7").Value = 45144.04206952546
$ws.Range("I7").ClearContents()
$ws.Range("J7").Value = '07/05/23 15:53'
$ws.Range("K7").Value = '07/05/23 15:53'
$ws.Range("M7").Value = '$2,100 as of 7/3/2023 1:33:47 PM'
$ws.Range("N7").Value = 2120

# Row 8
$ws.Range("A8").Value = 'L475090'
$ws.Range("C8").Value = 'S.B. 2'
$ws.Range("E8").Value = 2180
$ws.Range("H8").Value = 45115.04206952546
$ws.Range("J8").Value = '07/06/23 14:19'
$ws.Range("K8").Value = '07/06/23 13:47'
$ws.Range("M8").Value = '$2,260 as of 7/5/2023 11:40:47 AM'
$ws.Range("N8").Value = 2140

# Row 9
$ws.Range("A9").Value = 'L662336'
$ws.Range("C9").Value = 'SB#4 MONA MARKET'
$ws.Range("E9").Value = 2620
$ws.Range("H9").Value = 45119.04206952546
$ws.Range("J9").Value = '07/05/23 20:35'
$ws.Range("K9").Value = '07/05/23 20:35'
$ws.Range("L9").Value = 100
$ws.Range("M9").Value = '$2,620 as of 7/4/2023 6:15:15 PM'
$ws.Range("N9").Value = 2720

# Row 10
$ws.Range("A10").Value = 'L474761'
$ws.Range("C10").Value = 'BABS MARKET'
$ws.Range("E10").Value = 2700
$ws.Range("H10").Value = 45120.04206952546
$ws.Range("I10").ClearContents()
$ws.Range("J10").Value = '07/05/23 21:38'
$ws.Range("K10").Value = '07/05/23 21:38'
$ws.Range("L10").Value = 40
$ws.Range("M10").Value = '$2,700 as of 7/5/2023 9:56:34 AM'
$ws.Range("N10").Value = 2760

# Row 11
$ws.Range("A11").Value = 'LK561655'
$ws.Range("C11").Value = 'CRENSHAW CRAVOR #2'
$ws.Range("E11").Value = 2780
$ws.Range("H11").ClearContents()
$ws.Range("I11").Value = 'ATM Inactive greater than 48 minutes'
$ws.Range("J11").Value = '01/23/20 08:24'
$ws.Range("K11").Value = '01/23/20 08:24'
$ws.Range("M11").Value = '$2,780 as of 1/23/2020 6:24:32 AM'
$ws.Range("N11").Value = 2800

# Row 12
$ws.Range("A12").Value = 'L474792'
$ws.Range("C12").Value = 'NICK SHELL SERVICE'
$ws.Range("E12").Value = 2900
$ws.Range("H12").Value = 45146.04206952546
$ws.Range("I12").Value = 'ATM Inactive greater than 2000 minutes'
$ws.Range("J12").Value = '07/04/23 20:34'
$ws.Range("K12").Value = '07/04/23 20:34'
$ws.Range("M12").Value = '$2,900 as of 7/4/2023 6:34:42 PM'
$ws.Range("N12").Value = 2940

# Row 13
$ws.Range("A13").Value = 'L474817'
$ws.Range("C13").Value = 'SAFETY MARKET'
$ws.Range("E13").Value = 3000
$ws.Range("H13").Value = 45120.04206952546
$ws.Range("J13").Value = '07/06/23 14:23'
$ws.Range("K13").Value = '07/06/23 12:29'
$ws.Range("L13").Value = 120
$ws.Range("M13").Value = '$3,000 as of 7/5/2023 11:46:37 AM'
$ws.Range("N13").Value = 3000

# Row 14
$ws.Range("A14").Value = 'L688961'
$ws.Range("C14").Value = 'MONA MART'
$ws.Range("E14").Value = 3540
$ws.Range("H14").Value = 45216.04206952546
$ws.Range("J14").Value = '07/05/23 21:26'
$ws.Range("K14").Value = '07/05/23 21:26'
$ws.Range("M14").Value = '$3,540 as of 7/4/2023 1:55:10 PM'
$ws.Range("N14").Value = 3620

# Row 15
$ws.Range("A15").Value = 'L475182'
$ws.Range("C15").Value = 'LA ESQUINA DE ORO'
$ws.Range("E15").Value = 3800
$ws.Range("H15").ClearContents()
$ws.Range("I15").Value = 'ATM Inactive greater than 48 minutes'
$ws.Range("J15").Value = '09/16/20 16:57'
$ws.Range("K15").Value = '09/15/20 23:38'
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = '$3,800 as of 9/16/2020 1:28:00 PM'
$ws.Range("N15").Value = 3800

# Row 16
$ws.Range("E16").Value = 3820
$ws.Range("H16").Value = 45136.04206952546
$ws.Range("J16").Value = '07/05/23 15:23'
$ws.Range("K16").Value = '07/05/23 09:07'
$ws.Range("M16").Value = '$3,820 as of 7/5/2023 7:07:49 AM'
$ws.Range("N16").Value = 3820

# Row 17
$ws.Range("E17").Value = 3820
$ws.Range("H17").Value = 45127.04206952546
$ws.Range("I17").Value = 'ATM Inactive greater than 2000 minutes'
$ws.Range("J17").Value = '07/03/23 15:16'
$ws.Range("K17").Value = '07/03/23 15:16'
$ws.Range("M17").Value = '$3,820 as of 7/3/2023 1:16:30 PM'
$ws.Range("N17").Value = 3940

# Row 18
$ws.Range("A18").Value = 'L474746'
$ws.Range("C18").Value = 'ZACATES MARKET'
$ws.Range("E18").Value = 4320
$ws.Range("H18").Value = 45126.04206952546
$ws.Range("J18").Value = '07/05/23 16:19'
$ws.Range("K18").Value = '07/05/23 16:19'
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = '$4,320 as of 7/5/2023 11:18:44 AM'
$ws.Range("N18").Value = 4360

# Row 19
$ws.Range("A19").Value = 'L704741'
$ws.Range("C19").Value = 'W ADAMS COIN LAUNDRY'
$ws.Range("E19").Value = 4740
$ws.Range("H19").Value = 45122.04206952546
$ws.Range("J19").Value = '07/06/23 14:09'
$ws.Range("K19").Value = '07/06/23 11:36'
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = '$4,740 as of 7/5/2023 3:58:45 AM'
$ws.Range("N19").Value = 4700

# Row 20
$ws.Range("E20").Value = 5880
$ws.Range("H20").Value = 45119.04206952546
$ws.Range("J20").Value = '07/06/23 14:22'
$ws.Range("K20").Value = '07/06/23 13:10'
$ws.Range("M20").Value = '$5,880 as of 7/5/2023 10:41:48 AM'
$ws.Range("N20").Value = 5880

# Row 21
$ws.Range("E21").Value = 5960
$ws.Range("H21").Value = 45303.04206952546
$ws.Range("I21").Value = 'ATM Inactive greater than 2000 minutes'
$ws.Range("J21").Value = '07/04/23 22:16'
$ws.Range("K21").Value = '07/04/23 22:16'
$ws.Range("M21").Value = '$5,960 as of 7/4/2023 8:16:31 PM'
$ws.Range("N21").Value = 6000

# Row 22
$ws.Range("A22").Value = 'L688966'
$ws.Range("C22").Value = 'LACON MINI MART'
$ws.Range("E22").Value = 6520
$ws.Range("H22").Value = 45227.04206952546
$ws.Range("I22").Value = 'ATM Inactive greater than 2000 minutes'
$ws.Range("J22").Value = '07/04/23 16:10'
$ws.Range("K22").Value = '07/04/23 16:10'
$ws.Range("L22").Value = 20
$ws.Range("M22").Value = '$6,520 as of 7/4/2023 2:10:50 PM'
$ws.Range("N22").Value = 6620

# Row 23
$ws.Range("A23").Value = 'L678988'
$ws.Range("C23").Value = 'PAYELESS MARKET'
$ws.Range("E23").Value = 7000
$ws.Range("H23").Value = 45134.04206952546
$ws.Range("I23").ClearContents()
$ws.Range("J23").Value = '07/05/23 19:44'
$ws.Range("K23").Value = '07/05/23 19:44'
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = '$7,000 as of 7/3/2023 4:44:25 PM'
$ws.Range("N23").Value = 7160

# Row 24
$ws.Range("A24").Value = 'L697590'
$ws.Range("C24").Value = 'S B MARKET ST'
$ws.Range("E24").Value = 8780
$ws.Range("H24").Value = 45369.04206952546
$ws.Range("I24").Value = 'ATM Inactive greater than 2000 minutes'
$ws.Range("J24").Value = '06/29/23 11:36'
$ws.Range("K24").Value = '06/29/23 11:36'
$ws.Range("M24").Value = '$8,780 as of 6/29/2023 9:36:36 AM'
$ws.Range("N24").Value = 8800

# Row 25
$ws.Range("A25").Value = 'L697589'
$ws.Range("C25").Value = 'S B DISCOUNT MART'
$ws.Range("E25").Value = 9420
$ws.Range("H25").Value = 45119.04206952546
$ws.Range("J25").Value = '07/06/23 13:51'
$ws.Range("K25").Value = '07/06/23 13:51'
$ws.Range("L25").Value = 40
$ws.Range("M25").Value = '$9,440 as of 7/5/2023 11:48:53 AM'
$ws.Range("N25").Value = 9440

# Row 26
$ws.Range("A26").Value = 'LK923383'
$ws.Range("C26").Value = 'SAMYS PHONE CARDS'
$ws.Range("E26").Value = 10460
$ws.Range("H26").Value = 45127.04206952546
$ws.Range("J26").Value = '07/06/23 13:59'
$ws.Range("K26").Value = '07/06/23 13:59'
$ws.Range("L26").Value = 100
$ws.Range("M26").Value = '$10,540 as of 7/3/2023 7:51:00 PM'
$ws.Range("N26").Value = 10540

# Row 27
$ws.Range("E27").Value = 21380
$ws.Range("H27").Value = 45121.04206952546
$ws.Range("J27").Value = '07/06/23 13:57'
$ws.Range("K27").Value = '07/06/23 13:57'
$ws.Range("M27").Value = '$21,400 as of 7/5/2023 11:17:32 AM'
$ws.Range("N27").Value = 21400

# Row 28
$ws.Range("E28").Value = 121200
